$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new year column (P) by copying the formatting of the existing
# last year column (O) for the header border row, the year-label row,
# and the data row, then filling in the new values.
$ws.Range("O3").Copy() | Out-Null
$ws.Range("P3").PasteSpecial(-4122) | Out-Null

$ws.Range("O4").Copy() | Out-Null
$ws.Range("P4").PasteSpecial(-4122) | Out-Null
$ws.Range("P4").Value = 2022

$ws.Range("O5").Copy() | Out-Null
$ws.Range("P5").PasteSpecial(-4122) | Out-Null
$ws.Range("P5").Value = 2.6

# Revised figures for the existing 2019-2021 columns.
$ws.Range("M5").Value = 2.6
$ws.Range("N5").Value = 2.4
$ws.Range("O5").Value = 3.3

# Move the active selection to the (now-empty) top of the new column.
$ws.Range("P3").Select() | Out-Null
